# Apply cryptocurrency price/volume updates (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.592.53"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "3.437.17"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'578.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'147.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.480"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'7.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.61%  "
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "4.025.56"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "'28.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.80%  "
$ws.Range("D15").Value = "3.433.83"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "62.673.10"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'14.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "'9.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("D21").Value = "'387.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D23").Value = "'0.559"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "3.586.96"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("E26").Value = "  -3.17%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'7.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -9.43%  "
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("D37").Value = "'32.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "'170.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "3.471.85"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D43").Value = "'42.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("D45").Value = "'4.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("D47").Value = "2.556.04"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "'22.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.41%  "
$ws.Range("E51").Value = "  +0.06%  "
